# Apply cell value updates and insert missing AP column values
# as described by the commit "remove legacy data for consistency"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("adv-samples-hybrid")
$ws.Range("C5").Value = 58.50000262260437
$ws.Range("E5").Value = 58.50000262260437
$ws.Range("K5").Value = 58.25000286102295
$ws.Range("M5").Value = 59.00000214576721
$ws.Range("S5").Value = 0.5000013113021851
$ws.Range("U5").Value = 0.5000013113021851
$ws.Range("AA5").Value = 0.5000013113021851
$ws.Range("AC5").Value = 0.5000013113021851
$ws.Range("AI5").Value = 58.25000286102295
$ws.Range("C6").Value = 53.75000286102295
$ws.Range("E6").Value = 53.49999713897705
$ws.Range("K6").Value = 51.99999856948853
$ws.Range("M6").Value = 53.75000286102295
$ws.Range("S6").Value = 0.5
$ws.Range("U6").Value = 0.5
$ws.Range("AA6").Value = 0.5
$ws.Range("AC6").Value = 0.5
$ws.Range("AI6").Value = 53.75000286102295
$ws.Range("C7").Value = 49.00000286102295
$ws.Range("E7").Value = 48.49999737739563
$ws.Range("K7").Value = 48.74999713897705
$ws.Range("M7").Value = 48.49999737739563
$ws.Range("S7").Value = 0.5000014305114746
$ws.Range("U7").Value = 0.5000014305114746
$ws.Range("AA7").Value = 0.5000014305114746
$ws.Range("AC7").Value = 0.5000014305114746
$ws.Range("AI7").Value = 48.74999713897705
$ws.Range("C8").Value = 41.74999785423279
$ws.Range("E8").Value = 43.25000238418579
$ws.Range("K8").Value = 39.25000023841858
$ws.Range("M8").Value = 43.00000262260437
$ws.Range("S8").Value = 0.0000019073486328125
$ws.Range("U8").Value = 0.0000019073486328125
$ws.Range("AA8").Value = 0.0000019073486328125
$ws.Range("AC8").Value = 0.0000019073486328125
$ws.Range("AI8").Value = 41.24999833106995
$ws.Range("C9").Value = 36.00000262260436
$ws.Range("E9").Value = 36.00000262260436
$ws.Range("K9").Value = 36.00000262260436
$ws.Range("M9").Value = 36.25000238418577
$ws.Range("S9").Value = 0.2500009536743093
$ws.Range("U9").Value = 0.2500009536743093
$ws.Range("AA9").Value = 0.2500009536743093
$ws.Range("AC9").Value = 0.2500009536743093
$ws.Range("AI9").Value = 35.75000286102294
$ws.Range("C10").Value = 22.49999833106995
$ws.Range("E10").Value = 28.99999809265137
$ws.Range("K10").Value = 10.49999785423279
$ws.Range("M10").Value = 30.25000286102295
$ws.Range("S10").Value = 0.2500016689300537
$ws.Range("U10").Value = 0.2500016689300537
$ws.Range("AA10").Value = 0.2500016689300537
$ws.Range("AC10").Value = 0.2500016689300537
$ws.Range("AI10").Value = 30.25000286102295
$ws.Range("C11").Value = 6.000002384185791
$ws.Range("E11").Value = 20.75000023841858
$ws.Range("K11").Value = 23.74999737739563
$ws.Range("M11").Value = 24.75000238418579
$ws.Range("S11").Value = 0.5000016689300537
$ws.Range("U11").Value = 0.5000016689300537
$ws.Range("AA11").Value = 0.5000016689300537
$ws.Range("AC11").Value = 0.5000016689300537
$ws.Range("AI11").Value = 20.75000023841858
$ws.Range("C12").Value = 14.2500011920929
$ws.Range("E12").Value = 19.25000238418579
$ws.Range("K12").Value = 8.500000715255737
$ws.Range("M12").Value = 17.24999833106995
$ws.Range("S12").Value = 0.500002384185791
$ws.Range("U12").Value = 0.500002384185791
$ws.Range("AA12").Value = 0.500002384185791
$ws.Range("AC12").Value = 0.500002384185791
$ws.Range("AI12").Value = 14.00000143051147
$ws.Range("C13").Value = 0.2500021457672119
$ws.Range("E13").Value = 9.749999046325684
$ws.Range("K13").Value = 12.00000286102295
$ws.Range("M13").Value = 12.25000262260437
$ws.Range("S13").Value = 0.5000019073486328
$ws.Range("U13").Value = 0.5000019073486328
$ws.Range("AA13").Value = 0.5000019073486328
$ws.Range("AC13").Value = 0.5000019073486328
$ws.Range("AI13").Value = 9.499999284744264
$ws.Range("C14").Value = 3.249999761581421
$ws.Range("E14").Value = 5.249997854232788
$ws.Range("K14").Value = 0.2500026226043701
$ws.Range("M14").Value = 2.500000476837158
$ws.Range("S14").Value = 0.500002384185791
$ws.Range("U14").Value = 0.500002384185791
$ws.Range("AA14").Value = 0.500002384185791
$ws.Range("AC14").Value = 0.500002384185791
$ws.Range("AI14").Value = 7.250001907348633

$ws = $wb.Worksheets.Item("out-of-source-original")
$ws.Range("B5").Value = 0.01047491401511991
$ws.Range("D5").Value = 0.01047491401511991
$ws.Range("F5").Value = 0.01047491401511991
$ws.Range("H5").Value = 0.01454591125327909
$ws.Range("J5").Value = 0.02268194526511991
$ws.Range("L5").Value = 0.01861690849143827
$ws.Range("N5").Value = 0.01047491401511991
$ws.Range("P5").Value = 0.01454591125327909
$ws.Range("Z5").Value = 0.01454591125327909
$ws.Range("AB5").Value = 0.01047491401511991
$ws.Range("B6").Value = 0.01047491401511991
$ws.Range("D6").Value = 0.02268194526511991
$ws.Range("F6").Value = 0.01454591125327909
$ws.Range("H6").Value = 0.01861690849143827
$ws.Range("J6").Value = 0.04303097099143827
$ws.Range("L6").Value = 0.01454591125327909
$ws.Range("N6").Value = 0.01454591125327909
$ws.Range("P6").Value = 0.01861690849143827
$ws.Range("Z6").Value = 0.01047491401511991
$ws.Range("AB6").Value = 0.01047491401511991
$ws.Range("B7").Value = 0.01047491401511991
$ws.Range("D7").Value = 0.03082393974143827
$ws.Range("F7").Value = 0.01861690849143827
$ws.Range("H7").Value = 0.01861690849143827
$ws.Range("J7").Value = 0.01047491401511991
$ws.Range("L7").Value = 0.01454591125327909
$ws.Range("N7").Value = 0.01861690849143827
$ws.Range("P7").Value = 0.01861690849143827
$ws.Range("Z7").Value = 0.01047491401511991
$ws.Range("AB7").Value = 0.01047491401511991
$ws.Range("B8").Value = 0.01454591125327909
$ws.Range("D8").Value = 0.01047491401511991
$ws.Range("F8").Value = 0.01047491401511991
$ws.Range("H8").Value = 0.01047491401511991
$ws.Range("J8").Value = 0.01047491401511991
$ws.Range("L8").Value = 0.01047491401511991
$ws.Range("N8").Value = 0.01047491401511991
$ws.Range("P8").Value = 0.01047491401511991
$ws.Range("Z8").Value = 0.01047491401511991
$ws.Range("AB8").Value = 0.01047491401511991
$ws.Range("B9").Value = 0.02268194526511991
$ws.Range("D9").Value = 0.01861690849143827
$ws.Range("F9").Value = 0.01047491401511991
$ws.Range("H9").Value = 0.02268194526511991
$ws.Range("J9").Value = 0.03082393974143827
$ws.Range("L9").Value = 0.01454591125327909
$ws.Range("N9").Value = 0.01047491401511991
$ws.Range("P9").Value = 0.02268194526511991
$ws.Range("Z9").Value = 0.01454591125327909
$ws.Range("AB9").Value = 0.01047491401511991
$ws.Range("B10").Value = 0.01047491401511991
$ws.Range("D10").Value = 0.01047491401511991
$ws.Range("F10").Value = 0.01047491401511991
$ws.Range("H10").Value = 0.03488897651511991
$ws.Range("J10").Value = 0.01047491401511991
$ws.Range("L10").Value = 0.01047491401511991
$ws.Range("N10").Value = 0.01047491401511991
$ws.Range("P10").Value = 0.03488897651511991
$ws.Range("Z10").Value = 0.01454591125327909
$ws.Range("AB10").Value = 0.01047491401511991
$ws.Range("B11").Value = 0.02268194526511991
$ws.Range("D11").Value = 0.01454591125327909
$ws.Range("F11").Value = 0.01047491401511991
$ws.Range("H11").Value = 0.02268194526511991
$ws.Range("J11").Value = 0.01454591125327909
$ws.Range("L11").Value = 0.01047491401511991
$ws.Range("N11").Value = 0.01047491401511991
$ws.Range("P11").Value = 0.02268194526511991
$ws.Range("Z11").Value = 0.01047491401511991
$ws.Range("AB11").Value = 0.01047491401511991
$ws.Range("B12").Value = 0.01454591125327909
$ws.Range("D12").Value = 0.01047491401511991
$ws.Range("F12").Value = 0.01047491401511991
$ws.Range("H12").Value = 0.01454591125327909
$ws.Range("J12").Value = 0.01047491401511991
$ws.Range("L12").Value = 0.01047491401511991
$ws.Range("N12").Value = 0.01047491401511991
$ws.Range("P12").Value = 0.01454591125327909
$ws.Range("Z12").Value = 0.01454591125327909
$ws.Range("AB12").Value = 0.01047491401511991
$ws.Range("B13").Value = 0.01047491401511991
$ws.Range("D13").Value = 0.01454591125327909
$ws.Range("F13").Value = 0.01047491401511991
$ws.Range("H13").Value = 0.04303097099143827
$ws.Range("J13").Value = 0.01047491401511991
$ws.Range("L13").Value = 0.01047491401511991
$ws.Range("N13").Value = 0.01047491401511991
$ws.Range("P13").Value = 0.04303097099143827
$ws.Range("Z13").Value = 0.01047491401511991
$ws.Range("AB13").Value = 0.01047491401511991
$ws.Range("B14").Value = 0.01047491401511991
$ws.Range("D14").Value = 0.01047491401511991
$ws.Range("F14").Value = 0.01047491401511991
$ws.Range("H14").Value = 0.03082393974143827
$ws.Range("J14").Value = 0.01047491401511991
$ws.Range("L14").Value = 0.01454591125327909
$ws.Range("N14").Value = 0.01047491401511991
$ws.Range("P14").Value = 0.03082393974143827
$ws.Range("Z14").Value = 0.01047491401511991
$ws.Range("AB14").Value = 0.02675294250327909
$ws.Range("B21").Value = 0.01454591125327909
$ws.Range("D21").Value = 0.01861690849143827
$ws.Range("F21").Value = 0.01454591125327909
$ws.Range("H21").Value = 0.01047491401511991
$ws.Range("J21").Value = 0.01047491401511991
$ws.Range("L21").Value = 0.01861690849143827
$ws.Range("N21").Value = 0.01454591125327909
$ws.Range("P21").Value = 0.01047491401511991
$ws.Range("R21").Value = 0.01861690849143827
$ws.Range("T21").Value = 0.01861690849143827
$ws.Range("V21").Value = 0.01861690849143827
$ws.Range("X21").Value = 0.01047491401511991
$ws.Range("Z21").Value = 0.01861690849143827
$ws.Range("AB21").Value = 0.01861690849143827
$ws.Range("AD21").Value = 0.01861690849143827
$ws.Range("AF21").Value = 0.01047491401511991
$ws.Range("AH21").Value = 0.01861690849143827
$ws.Range("AJ21").Value = 0.01047491401511991
$ws.Range("AL21").Value = 0.01861690849143827
$ws.Range("AN21").Value = 0.01454591125327909
$ws.Range("B22").Value = 0.01454591125327909
$ws.Range("D22").Value = 0.01047491401511991
$ws.Range("F22").Value = 0.01047491401511991
$ws.Range("H22").Value = 0.01047491401511991
$ws.Range("J22").Value = 0.01454591125327909
$ws.Range("L22").Value = 0.01861690849143827
$ws.Range("N22").Value = 0.01047491401511991
$ws.Range("P22").Value = 0.01047491401511991
$ws.Range("R22").Value = 0.01861690849143827
$ws.Range("T22").Value = 0.01861690849143827
$ws.Range("V22").Value = 0.02268194526511991
$ws.Range("X22").Value = 0.01047491401511991
$ws.Range("Z22").Value = 0.01861690849143827
$ws.Range("AB22").Value = 0.01861690849143827
$ws.Range("AD22").Value = 0.02268194526511991
$ws.Range("AF22").Value = 0.01047491401511991
$ws.Range("AH22").Value = 0.01454591125327909
$ws.Range("AJ22").Value = 0.01861690849143827
$ws.Range("AL22").Value = 0.02268194526511991
$ws.Range("AN22").Value = 0.01454591125327909
$ws.Range("B23").Value = 0.01047491401511991
$ws.Range("D23").Value = 0.01861690849143827
$ws.Range("F23").Value = 0.01047491401511991
$ws.Range("H23").Value = 0.01047491401511991
$ws.Range("J23").Value = 0.01861690849143827
$ws.Range("L23").Value = 0.01454591125327909
$ws.Range("N23").Value = 0.01047491401511991
$ws.Range("P23").Value = 0.01047491401511991
$ws.Range("R23").Value = 0.01861690849143827
$ws.Range("T23").Value = 0.02268194526511991
$ws.Range("V23").Value = 0.02268194526511991
$ws.Range("X23").Value = 0.01454591125327909
$ws.Range("Z23").Value = 0.01861690849143827
$ws.Range("AB23").Value = 0.02268194526511991
$ws.Range("AD23").Value = 0.02268194526511991
$ws.Range("AF23").Value = 0.01454591125327909
$ws.Range("AH23").Value = 0.01861690849143827
$ws.Range("AJ23").Value = 0.01454591125327909
$ws.Range("AL23").Value = 0.02675294250327909
$ws.Range("AN23").Value = 0.02268194526511991
$ws.Range("B24").Value = 0.03488897651511991
$ws.Range("D24").Value = 0.01454591125327909
$ws.Range("F24").Value = 0.01454591125327909
$ws.Range("H24").Value = 0.01454591125327909
$ws.Range("J24").Value = 0.02675294250327909
$ws.Range("L24").Value = 0.02268194526511991
$ws.Range("N24").Value = 0.01454591125327909
$ws.Range("P24").Value = 0.01454591125327909
$ws.Range("R24").Value = 0.01454591125327909
$ws.Range("T24").Value = 0.02268194526511991
$ws.Range("V24").Value = 0.02268194526511991
$ws.Range("X24").Value = 0.01861690849143827
$ws.Range("Z24").Value = 0.01454591125327909
$ws.Range("AB24").Value = 0.02268194526511991
$ws.Range("AD24").Value = 0.02268194526511991
$ws.Range("AF24").Value = 0.01861690849143827
$ws.Range("AH24").Value = 0.02675294250327909
$ws.Range("AJ24").Value = 0.01454591125327909
$ws.Range("AL24").Value = 0.02675294250327909
$ws.Range("AN24").Value = 0.02675294250327909
$ws.Range("B25").Value = 0.01047491401511991
$ws.Range("D25").Value = 0.01047491401511991
$ws.Range("F25").Value = 0.01454591125327909
$ws.Range("H25").Value = 0.01047491401511991
$ws.Range("J25").Value = 0.01047491401511991
$ws.Range("L25").Value = 0.01047491401511991
$ws.Range("N25").Value = 0.01454591125327909
$ws.Range("P25").Value = 0.01047491401511991
$ws.Range("R25").Value = 0.02268194526511991
$ws.Range("T25").Value = 0.02268194526511991
$ws.Range("V25").Value = 0.02268194526511991
$ws.Range("X25").Value = 0.01861690849143827
$ws.Range("Z25").Value = 0.02268194526511991
$ws.Range("AB25").Value = 0.02268194526511991
$ws.Range("AD25").Value = 0.02268194526511991
$ws.Range("AF25").Value = 0.01861690849143827
$ws.Range("AH25").Value = 0.01861690849143827
$ws.Range("AJ25").Value = 0.01454591125327909
$ws.Range("AL25").Value = 0.01454591125327909
$ws.Range("AN25").Value = 0.01861690849143827
$ws.Range("B26").Value = 0.01861690849143827
$ws.Range("D26").Value = 0.01454591125327909
$ws.Range("F26").Value = 0.01454591125327909
$ws.Range("H26").Value = 0.01047491401511991
$ws.Range("J26").Value = 0.01454591125327909
$ws.Range("L26").Value = 0.01047491401511991
$ws.Range("N26").Value = 0.01454591125327909
$ws.Range("P26").Value = 0.01047491401511991
$ws.Range("R26").Value = 0.02268194526511991
$ws.Range("T26").Value = 0.02268194526511991
$ws.Range("V26").Value = 0.03488897651511991
$ws.Range("X26").Value = 0.02268194526511991
$ws.Range("Z26").Value = 0.01861690849143827
$ws.Range("AB26").Value = 0.02268194526511991
$ws.Range("AD26").Value = 0.03488897651511991
$ws.Range("AF26").Value = 0.02268194526511991
$ws.Range("AH26").Value = 0.04709600776511991
$ws.Range("AJ26").Value = 0.02268194526511991
$ws.Range("AL26").Value = 0.03488897651511991
$ws.Range("AN26").Value = 0.03082393974143827
$ws.Range("B27").Value = 0.01047491401511991
$ws.Range("D27").Value = 0.01047491401511991
$ws.Range("F27").Value = 0.01047491401511991
$ws.Range("H27").Value = 0.02268194526511991
$ws.Range("J27").Value = 0.03488897651511991
$ws.Range("L27").Value = 0.01861690849143827
$ws.Range("N27").Value = 0.01047491401511991
$ws.Range("P27").Value = 0.02268194526511991
$ws.Range("R27").Value = 0.02268194526511991
$ws.Range("T27").Value = 0.02675294250327909
$ws.Range("V27").Value = 0.03082393974143827
$ws.Range("X27").Value = 0.02268194526511991
$ws.Range("Z27").Value = 0.02268194526511991
$ws.Range("AB27").Value = 0.02675294250327909
$ws.Range("AD27").Value = 0.03082393974143827
$ws.Range("AF27").Value = 0.02268194526511991
$ws.Range("AH27").Value = 0.01861690849143827
$ws.Range("AJ27").Value = 0.04303097099143827
$ws.Range("AL27").Value = 0.03488897651511991
$ws.Range("AN27").Value = 0.03082393974143827
$ws.Range("B28").Value = 0.01861690849143827
$ws.Range("D28").Value = 0.01047491401511991
$ws.Range("F28").Value = 0.02268194526511991
$ws.Range("H28").Value = 0.01047491401511991
$ws.Range("J28").Value = 0.01047491401511991
$ws.Range("L28").Value = 0.01047491401511991
$ws.Range("N28").Value = 0.02268194526511991
$ws.Range("P28").Value = 0.01047491401511991
$ws.Range("R28").Value = 0.01861690849143827
$ws.Range("T28").Value = 0.01861690849143827
$ws.Range("V28").Value = 0.03488897651511991
$ws.Range("X28").Value = 0.06337403625327909
$ws.Range("Z28").Value = 0.01861690849143827
$ws.Range("AB28").Value = 0.01861690849143827
$ws.Range("AD28").Value = 0.03488897651511991
$ws.Range("AF28").Value = 0.06337403625327909
$ws.Range("AH28").Value = 0.01454591125327909
$ws.Range("AJ28").Value = 0.03082393974143827
$ws.Range("AL28").Value = 0.03082393974143827
$ws.Range("AN28").Value = 0.02268194526511991
$ws.Range("B29").Value = 0.01454591125327909
$ws.Range("D29").Value = 0.01047491401511991
$ws.Range("F29").Value = 0.01047491401511991
$ws.Range("H29").Value = 0.01047491401511991
$ws.Range("J29").Value = 0.01454591125327909
$ws.Range("L29").Value = 0.01047491401511991
$ws.Range("N29").Value = 0.01047491401511991
$ws.Range("P29").Value = 0.01047491401511991
$ws.Range("R29").Value = 0.02268194526511991
$ws.Range("T29").Value = 0.01861690849143827
$ws.Range("V29").Value = 0.05116700500327909
$ws.Range("X29").Value = 0.05523800224143827
$ws.Range("Z29").Value = 0.02268194526511991
$ws.Range("AB29").Value = 0.01861690849143827
$ws.Range("AD29").Value = 0.05116700500327909
$ws.Range("AF29").Value = 0.05523800224143827
$ws.Range("AH29").Value = 0.01454591125327909
$ws.Range("AJ29").Value = 0.01047491401511991
$ws.Range("AL29").Value = 0.01454591125327909
$ws.Range("AN29").Value = 0.01861690849143827
$ws.Range("B30").Value = 0.01454591125327909
$ws.Range("D30").Value = 0.02268194526511991
$ws.Range("F30").Value = 0.01047491401511991
$ws.Range("H30").Value = 0.01861690849143827
$ws.Range("J30").Value = 0.01047491401511991
$ws.Range("L30").Value = 0.01047491401511991
$ws.Range("N30").Value = 0.01047491401511991
$ws.Range("P30").Value = 0.01861690849143827
$ws.Range("R30").Value = 0.04303097099143827
$ws.Range("T30").Value = 0.03895997375327909
$ws.Range("V30").Value = 0.05523800224143827
$ws.Range("X30").Value = 0.06744503349143828
$ws.Range("Z30").Value = 0.03082393974143827
$ws.Range("AB30").Value = 0.03488897651511991
$ws.Range("AD30").Value = 0.05523800224143827
$ws.Range("AF30").Value = 0.06744503349143828
$ws.Range("AH30").Value = 0.01047491401511991
$ws.Range("AJ30").Value = 0.01454591125327909
$ws.Range("AL30").Value = 0.01454591125327909
$ws.Range("AN30").Value = 0.02268194526511991

$ws = $wb.Worksheets.Item("out-of-source-hybrid")
$ws.Range("AP37").Value = 2.624997568130496
$ws.Range("AP38").Value = 2.700000190734869
$ws.Range("AP39").Value = 2.799998426437384
$ws.Range("AP40").Value = 2.825001049041759
$ws.Range("AP41").Value = 2.724999761581429
$ws.Range("AP42").Value = 2.62500185966492
$ws.Range("AP43").Value = 2.174997472763067
$ws.Range("AP44").Value = 1.699997758865351
$ws.Range("AP45").Value = 1.325000381469735
$ws.Range("AP46").Value = 0.57499809265137
$ws.Range("AP47").Value = -0.000002384185791015625
